$wb = $excel.ActiveWorkbook

# Sheet: 展览 (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F6").Value = 7953
$ws1.Range("F9").Value = 23
$ws1.Range("F10").Value = 1911
$ws1.Range("F11").Value = 435
$ws1.Range("F12").Value = 188
$ws1.Range("F13").Value = 1792
$ws1.Range("F15").Value = 1122
$ws1.Range("F18").Value = 1104
$ws1.Range("F19").Value = 8607
$ws1.Range("F20").Value = 219
$ws1.Range("F25").Value = 1044
$ws1.Range("F26").Value = 584
$ws1.Range("F27").Value = 1209
$ws1.Range("F28").Value = 1082
$ws1.Range("F29").Value = 611
$ws1.Range("F33").Value = 126
$ws1.Range("F34").Value = 1069
$ws1.Range("F35").Value = 484
$ws1.Range("F36").Value = 384
$ws1.Range("F37").Value = 3629
$ws1.Range("F40").Value = 3
$ws1.Range("F41").Value = 532
$ws1.Range("F43").Value = 6
$ws1.Range("F44").Value = 733
$ws1.Range("F46").Value = 121
$ws1.Range("F48").Value = 37
$ws1.Range("F49").Value = 5

# Sheet: 演出 (Performance)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 360
$ws2.Range("F11").Value = 37
$ws2.Range("F12").Value = 37
$ws2.Range("F21").Value = 55
$ws2.Range("F25").Value = 7033
$ws2.Range("F41").Value = 1

# Sheet: 本地生活 (Local Life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F4").Value = 2113
$ws3.Range("F5").Value = 1423
$ws3.Range("F8").Value = 2266
$ws3.Range("F9").Value = 9093
$ws3.Range("F10").Value = 1376

# Sheet: 全部类型 (All Types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 2113
$ws4.Range("F5").Value = 7953
$ws4.Range("F7").Value = 2266
$ws4.Range("F8").Value = 1376
$ws4.Range("F10").Value = 435
$ws4.Range("F11").Value = 188
$ws4.Range("F12").Value = 1792
$ws4.Range("F14").Value = 1122
$ws4.Range("F17").Value = 1104
$ws4.Range("F18").Value = 8607
$ws4.Range("F19").Value = 219
$ws4.Range("F23").Value = 1044
$ws4.Range("F24").Value = 584
$ws4.Range("F25").Value = 1209
$ws4.Range("F26").Value = 1082
$ws4.Range("F27").Value = 611
$ws4.Range("F30").Value = 37
$ws4.Range("F32").Value = 126
$ws4.Range("F33").Value = 1069
$ws4.Range("F34").Value = 484
$ws4.Range("F37").Value = 3629
$ws4.Range("F39").Value = 532
$ws4.Range("F41").Value = 6
$ws4.Range("F42").Value = 733
$ws4.Range("F45").Value = 121
$ws4.Range("F48").Value = 37
